# Applies the 2024-09-21 "android.xlsx" edit:
#  - A new September SMS entry ("service axis" @ 2024-09-21 19:00:35) was
#    recorded at the top of the September_Details/September_Date columns
#    (R/S) on the "2024" sheet, pushing every existing R/S pair down by one
#    row (row N -> row N+1) for rows 3..185.
#  - Because column A ("Group") entries live interleaved in that same
#    column range, they cascade down by the same one row for rows 3..185
#    (row 2's "Mobile" label is untouched).
#  - The last Group label ("Broadband", previously row 185) now lands in a
#    newly created row 186, which grows the sheet's used range from
#    A1:Y185 to A1:Y186.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("2024")

$firstShiftRow = 3
$lastShiftRow = 185

# Snapshot the current (pre-edit) column A / R / S values for the rows that
# are about to move, reading from the bottom up is not required here since
# we buffer everything in arrays first, then write afterwards.
$colA = @{}
$colR = @{}
$colS = @{}
for ($r = $firstShiftRow; $r -le $lastShiftRow; $r++) {
    $colA[$r] = $ws.Cells.Item($r, 1).Value2
    $colR[$r] = $ws.Cells.Item($r, 18).Value2
    $colS[$r] = $ws.Cells.Item($r, 19).Value2
}

# Write the snapshot back out shifted down by one row (process bottom-up so
# we never clobber a source row before it has been read -- not strictly
# needed since we buffered above, but keeps the intent obvious/safe).
# Cells whose destination already holds the same value are left untouched
# so we don't needlessly disturb cells that aren't actually changing.
for ($r = $lastShiftRow; $r -ge $firstShiftRow; $r--) {
    $target = $r + 1

    if ($ws.Cells.Item($target, 1).Value2 -ne $colA[$r]) {
        $ws.Cells.Item($target, 1).Value = $colA[$r]
    }
    if ($ws.Cells.Item($target, 18).Value2 -ne $colR[$r]) {
        $ws.Cells.Item($target, 18).Value = $colR[$r]
    }
    if ($ws.Cells.Item($target, 19).Value2 -ne $colS[$r]) {
        $ws.Cells.Item($target, 19).Value = $colS[$r]
    }
}

# The brand-new SMS entry occupies the freed-up row 3; it has no Group
# label of its own (column A stays blank there, same as before the edit).
$ws.Cells.Item($firstShiftRow, 1).Value = ""
$ws.Cells.Item($firstShiftRow, 18).Value = "service axis"
$ws.Cells.Item($firstShiftRow, 19).Value = "2024-09-21 19:00:35"
